$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 155
$ws.Range("J2").Value = 612
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 143
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 107
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 75
$ws.Range("T2").Value = 116
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 991
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 977
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 3
